# Correct the dish_type data in column E: replace the abbreviated
# "NV"/"V" labels with descriptive Main/Starter categories.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = "Non-Veg Main"
$ws.Range("E3").Value  = "Non-Veg Main"
$ws.Range("E4").Value  = "Non-Veg Main"
$ws.Range("E5").Value  = "Veg Main"
$ws.Range("E6").Value  = "Veg Main"
$ws.Range("E7").Value  = "Veg Main"
$ws.Range("E8").Value  = "Non-Veg Starter"
$ws.Range("E9").Value  = "Non-Veg Starter"
$ws.Range("E10").Value = "Non-Veg Main"
$ws.Range("E11").Value = "Veg Starter"
$ws.Range("E12").Value = "Non-Veg Starter"
$ws.Range("E13").Value = "Veg Main"
$ws.Range("E14").Value = "Veg Main"
$ws.Range("E15").Value = "Non-Veg Main"
$ws.Range("E16").Value = "Non-Veg Main"
$ws.Range("E17").Value = "Non-Veg Main"
$ws.Range("E18").Value = "Veg Main"
$ws.Range("E19").Value = "Veg Main"
$ws.Range("E20").Value = "Veg Main"
$ws.Range("E21").Value = "Veg Main"

# Update the active selection to match the latest edit location.
$ws.Range("E22").Select()
